$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume snapshot data. Two pairs of adjacent rows
# also swapped rank position (Litecoin/ShibaInu at 15-16, and
# RocketPoolETH/FraxShare at 42-43).

$ws.Cells.Item(2, 4).Value = '29.404.57'
$ws.Cells.Item(2, 5).Value = '  -0.52%  '

$ws.Cells.Item(3, 4).Value = '1.849.74'
$ws.Cells.Item(3, 5).Value = '  -0.06%  '

$ws.Cells.Item(4, 4).Value = '''0.9994'
$ws.Cells.Item(4, 5).Value = '  -0.16%  '

$ws.Cells.Item(5, 4).Value = '''241.60'
$ws.Cells.Item(5, 5).Value = '  -0.75%  '

$ws.Cells.Item(6, 4).Value = '''0.6276'
$ws.Cells.Item(6, 5).Value = '  -3.82%  '

$ws.Cells.Item(7, 4).Value = '''1.000'
$ws.Cells.Item(7, 5).Value = '  -0.10%  '

$ws.Cells.Item(8, 4).Value = '''0.07644'
$ws.Cells.Item(8, 5).Value = '  +2.15%  '

$ws.Cells.Item(9, 4).Value = '''0.2973'
$ws.Cells.Item(9, 5).Value = '  +0.05%  '

$ws.Cells.Item(10, 4).Value = '''24.50'
$ws.Cells.Item(10, 5).Value = '  +0.16%  '

$ws.Cells.Item(11, 4).Value = '2.040.14'
$ws.Cells.Item(11, 5).Value = '  +10.04%  '

$ws.Cells.Item(12, 4).Value = '''0.07717'
$ws.Cells.Item(12, 5).Value = '  +1.00%  '

$ws.Cells.Item(13, 4).Value = '''5.003'
$ws.Cells.Item(13, 5).Value = '  -0.72%  '

$ws.Cells.Item(14, 4).Value = '''0.6890'
$ws.Cells.Item(14, 5).Value = '  +0.54%  '

$ws.Cells.Item(15, 2).Value = 'Litecoin'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(15, 4).Value = '''83.06'
$ws.Cells.Item(15, 5).Value = '  -0.46%  '

$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(16, 4).Value = '''0.000009971'
$ws.Cells.Item(16, 5).Value = '  +4.13%  '

$ws.Cells.Item(17, 4).Value = '2.220.81'
$ws.Cells.Item(17, 5).Value = '  +5.16%  '

$ws.Cells.Item(18, 4).Value = '''6.190'
$ws.Cells.Item(18, 5).Value = '  +1.07%  '

$ws.Cells.Item(19, 4).Value = '29.612.98'
$ws.Cells.Item(19, 5).Value = '  +0.12%  '

$ws.Cells.Item(20, 4).Value = '''232.31'
$ws.Cells.Item(20, 5).Value = '  -1.86%  '

$ws.Cells.Item(21, 4).Value = '''12.55'
$ws.Cells.Item(21, 5).Value = '  -0.40%  '

$ws.Cells.Item(22, 4).Value = '''1.001'
$ws.Cells.Item(22, 5).Value = '  +0.00%  '

$ws.Cells.Item(23, 4).Value = '''7.663'
$ws.Cells.Item(23, 5).Value = '  -0.43%  '

$ws.Cells.Item(24, 4).Value = '''1.000'
$ws.Cells.Item(24, 5).Value = '  -0.26%  '

$ws.Cells.Item(25, 4).Value = '''155.03'
$ws.Cells.Item(25, 5).Value = '  -1.68%  '

$ws.Cells.Item(26, 4).Value = '''0.1388'
$ws.Cells.Item(26, 5).Value = '  -2.48%  '

$ws.Cells.Item(27, 4).Value = '''8.471'
$ws.Cells.Item(27, 5).Value = '  -0.51%  '

$ws.Cells.Item(28, 5).Value = '  -0.79%  '

$ws.Cells.Item(29, 4).Value = '''1.474'
$ws.Cells.Item(29, 5).Value = '  -1.13%  '

$ws.Cells.Item(30, 4).Value = '''0.05790'
$ws.Cells.Item(30, 5).Value = '  -4.03%  '

$ws.Cells.Item(31, 5).Value = '  -0.19%  '

$ws.Cells.Item(32, 5).Value = '  -0.05%  '

$ws.Cells.Item(33, 4).Value = '''4.019'
$ws.Cells.Item(33, 5).Value = '  -1.37%  '

$ws.Cells.Item(34, 4).Value = '''1.877'
$ws.Cells.Item(34, 5).Value = '  +0.45%  '

$ws.Cells.Item(35, 4).Value = '''1.160'
$ws.Cells.Item(35, 5).Value = '  -2.12%  '

$ws.Cells.Item(36, 4).Value = '''0.7209'

$ws.Cells.Item(37, 4).Value = '''2.584'
$ws.Cells.Item(37, 5).Value = '  -0.60%  '

$ws.Cells.Item(38, 4).Value = '1.256.66'
$ws.Cells.Item(38, 5).Value = '  +4.66%  '

$ws.Cells.Item(39, 4).Value = '''2.796'
$ws.Cells.Item(39, 5).Value = '  -0.16%  '

$ws.Cells.Item(40, 4).Value = '''0.01806'
$ws.Cells.Item(40, 5).Value = '  +1.16%  '

$ws.Cells.Item(41, 4).Value = '''0.9092'
$ws.Cells.Item(41, 5).Value = '  -0.32%  '

$ws.Cells.Item(42, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(42, 4).Value = '2.164.06'
$ws.Cells.Item(42, 5).Value = '  +7.27%  '

$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).Value = '''6.087'
$ws.Cells.Item(43, 5).Value = '  -3.19%  '

$ws.Cells.Item(44, 4).Value = '''0.9993'
$ws.Cells.Item(44, 5).Value = '  -0.21%  '

$ws.Cells.Item(45, 4).Value = '''67.58'
$ws.Cells.Item(45, 5).Value = '  +1.66%  '

$ws.Cells.Item(46, 4).Value = '''101.71'
$ws.Cells.Item(46, 5).Value = '  +0.30%  '

$ws.Cells.Item(47, 4).Value = '''7.299'
$ws.Cells.Item(47, 5).Value = '  -0.54%  '

$ws.Cells.Item(48, 5).Value = '  -3.28%  '

$ws.Cells.Item(49, 4).Value = '''9.165'
$ws.Cells.Item(49, 5).Value = '  +0.57%  '

$ws.Cells.Item(50, 4).Value = '''0.4032'
$ws.Cells.Item(50, 5).Value = '  -0.44%  '

$ws.Cells.Item(51, 4).Value = '''1.698'
$ws.Cells.Item(51, 5).Value = '  +2.20%  '
